# Generate Report for Handoff
# Adds two new handed-off files (87591eed-... and b7f00d72-...) as new rows
# to the Overview sheet, the zh-cn sheet and the de-de sheet, and grows the
# three tables / sheet dimensions accordingly.

$wb = $excel.ActiveWorkbook

$HYPERLINK_COLOR = 15570276   # BGR encoding of RGB(0x64,0x95,0xED) -> matches the workbook's HyperLink font color
$DATEFMT = "yyyy-mm-dd HH:mm:ss"

function Set-HandoffHyperlink($ws, $range, $address, $displayText) {
    $ws.Hyperlinks.Add($range, $address, "", "", $displayText) | Out-Null
    $range.Font.Underline = 2
    $range.Font.Color = $HYPERLINK_COLOR
}

function Set-DateText($range, $text) {
    $range.Value = $text
    $range.NumberFormat = $DATEFMT
}

# "True"/"False" look like booleans to the COM layer's type-inference, so they
# get stored as real xlsx booleans (t="b") unless forced to text with a
# leading apostrophe, exactly like typing '\''False into Excel would.
function Set-BooleanText($range, $text) {
    $range.Value = "'" + $text
}

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) - columns: A File Name, B Path And Name,
# C Extension, D Publish URL, E zh-cn, F de-de, G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 4 - 87591eed-a857-4183-afeb-ded04c3df193.md
$wsOverview.Range("A4").Value = "87591eed-a857-4183-afeb-ded04c3df193.md"
Set-HandoffHyperlink $wsOverview $wsOverview.Range("B4") `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87591eed-a857-4183-afeb-ded04c3df193/e2e/87591eed-a857-4183-afeb-ded04c3df193.md" `
    "e2e\87591eed-a857-4183-afeb-ded04c3df193.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
Set-DateText $wsOverview.Range("G4") "2016-09-06 15:25:30"

# Row 5 - b7f00d72-5925-4a8b-801c-39e1ccef97a8.md
$wsOverview.Range("A5").Value = "b7f00d72-5925-4a8b-801c-39e1ccef97a8.md"
Set-HandoffHyperlink $wsOverview $wsOverview.Range("B5") `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7f00d72-5925-4a8b-801c-39e1ccef97a8/e2e/b7f00d72-5925-4a8b-801c-39e1ccef97a8.md" `
    "e2e\b7f00d72-5925-4a8b-801c-39e1ccef97a8.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
Set-DateText $wsOverview.Range("G5") "2016-09-06 15:25:30"

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G5")) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) - 16 columns A..P
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 4
Set-HandoffHyperlink $wsZhCn $wsZhCn.Range("A4") `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87591eed-a857-4183-afeb-ded04c3df193/e2e/87591eed-a857-4183-afeb-ded04c3df193.md" `
    "87591eed-a857-4183-afeb-ded04c3df193.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
Set-BooleanText $wsZhCn.Range("F4") "False"
$wsZhCn.Range("G4").Value = "87591eed-a857-4183-afeb-ded04c3df193.b614a63b8e57bd40671afb5b30e29c144e084624.zh-cn.xlf"
Set-DateText $wsZhCn.Range("H4") "2016-09-06 15:25:13"
$wsZhCn.Range("I4").Value = ""
$wsZhCn.Range("J4").Value = ""
Set-DateText $wsZhCn.Range("K4") "0001-01-01 00:00:00"
$wsZhCn.Range("L4").Value = ""
Set-BooleanText $wsZhCn.Range("M4") "True"
$wsZhCn.Range("N4").Value = ""
Set-BooleanText $wsZhCn.Range("O4") "False"
$wsZhCn.Range("P4").Value = ""

# Row 5
Set-HandoffHyperlink $wsZhCn $wsZhCn.Range("A5") `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7f00d72-5925-4a8b-801c-39e1ccef97a8/e2e/b7f00d72-5925-4a8b-801c-39e1ccef97a8.md" `
    "b7f00d72-5925-4a8b-801c-39e1ccef97a8.md"
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
Set-BooleanText $wsZhCn.Range("F5") "False"
$wsZhCn.Range("G5").Value = "b7f00d72-5925-4a8b-801c-39e1ccef97a8.604b4b4686bacecbb74355f6ce8cf80a5050ecb3.zh-cn.xlf"
Set-DateText $wsZhCn.Range("H5") "2016-09-06 15:25:13"
$wsZhCn.Range("I5").Value = ""
$wsZhCn.Range("J5").Value = ""
Set-DateText $wsZhCn.Range("K5") "0001-01-01 00:00:00"
$wsZhCn.Range("L5").Value = ""
Set-BooleanText $wsZhCn.Range("M5") "True"
$wsZhCn.Range("N5").Value = ""
Set-BooleanText $wsZhCn.Range("O5") "False"
$wsZhCn.Range("P5").Value = ""

$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:P5")) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) - 16 columns A..P
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 4
Set-HandoffHyperlink $wsDeDe $wsDeDe.Range("A4") `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87591eed-a857-4183-afeb-ded04c3df193/e2e/87591eed-a857-4183-afeb-ded04c3df193.md" `
    "87591eed-a857-4183-afeb-ded04c3df193.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
Set-BooleanText $wsDeDe.Range("F4") "False"
$wsDeDe.Range("G4").Value = "87591eed-a857-4183-afeb-ded04c3df193.b614a63b8e57bd40671afb5b30e29c144e084624.de-de.xlf"
Set-DateText $wsDeDe.Range("H4") "2016-09-06 15:25:30"
$wsDeDe.Range("I4").Value = ""
$wsDeDe.Range("J4").Value = ""
Set-DateText $wsDeDe.Range("K4") "0001-01-01 00:00:00"
$wsDeDe.Range("L4").Value = ""
Set-BooleanText $wsDeDe.Range("M4") "True"
$wsDeDe.Range("N4").Value = ""
Set-BooleanText $wsDeDe.Range("O4") "False"
$wsDeDe.Range("P4").Value = ""

# Row 5
Set-HandoffHyperlink $wsDeDe $wsDeDe.Range("A5") `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7f00d72-5925-4a8b-801c-39e1ccef97a8/e2e/b7f00d72-5925-4a8b-801c-39e1ccef97a8.md" `
    "b7f00d72-5925-4a8b-801c-39e1ccef97a8.md"
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
Set-BooleanText $wsDeDe.Range("F5") "False"
$wsDeDe.Range("G5").Value = "b7f00d72-5925-4a8b-801c-39e1ccef97a8.604b4b4686bacecbb74355f6ce8cf80a5050ecb3.de-de.xlf"
Set-DateText $wsDeDe.Range("H5") "2016-09-06 15:25:30"
$wsDeDe.Range("I5").Value = ""
$wsDeDe.Range("J5").Value = ""
Set-DateText $wsDeDe.Range("K5") "0001-01-01 00:00:00"
$wsDeDe.Range("L5").Value = ""
Set-BooleanText $wsDeDe.Range("M5") "True"
$wsDeDe.Range("N5").Value = ""
Set-BooleanText $wsDeDe.Range("O5") "False"
$wsDeDe.Range("P5").Value = ""

$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:P5")) | Out-Null
